# Adding network health monitoring
# Insert "Address" and "Phone" columns into the Examinees header row,
# and mirror the full header row onto the (previously empty) Examiners sheet.

$wb = $excel.ActiveWorkbook

$wsExaminees = $wb.Worksheets.Item("Examinees")
$wsExaminers = $wb.Worksheets.Item("Examiners")

# --- Examinees sheet: insert two new header columns ---------------------
# Current header (A..I): ID, First Name, Last Name, KKW ID, NOC, DOB, Email, Group, New Rank
# Target header (A..K): ID, First Name, Last Name, KKW ID, NOC, DOB, Address, Email, Phone, Group, New Rank

# Insert the two new (blank) columns first ...
$wsExaminees.Range("G1").EntireColumn.Insert()
$wsExaminees.Range("I1").EntireColumn.Insert()

# ... then fill them in, Phone (I) before Address (G)
$wsExaminees.Range("I1").Value = "Phone"
$wsExaminees.Range("G1").Value = "Address"

# Widen column A
$wsExaminees.Columns.Item(1).ColumnWidth = 9.9167

# Update selection on Examinees
[void]$wsExaminees.Range("G3").Select()

# --- Examiners sheet: populate header row matching Examinees (minus New Rank) ---
$headers = @("ID", "First Name", "Last Name", "KKW ID", "NOC", "DOB", "Address", "Email", "Phone", "Group")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $wsExaminers.Cells.Item(1, $col).Value = $headers[$i]
}

[void]$wsExaminers.Range("G3").Select()

# --- Make Examiners the active sheet/tab ---------------------------------
$wsExaminers.Activate()
